$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the new value for C2 (3' Barcode Sequence for Sample1)
$ws.Range("C2").Value = "NNGACNN"

# Update the selection to C2 to match the saved workbook view state
$ws.Range("C2").Select()
